# Regenerate handback-status report: refresh processing timestamps for the
# "15c4020e..." and "e0e78c73..." entries (rows 2 and 5 in each sheet), and
# flip their Priority from "ht" to "mt".

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-13 04:18:49"
$wsOverview.Range("G5").Value = "2016-08-13 04:18:49"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-13 04:18:41"
$wsZhCn.Range("H5").Value = "2016-08-13 04:18:41"
$wsZhCn.Range("K2").Value = "2016-08-13 04:19:13"
$wsZhCn.Range("K5").Value = "2016-08-13 04:19:13"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-13 04:18:49"
$wsDeDe.Range("H5").Value = "2016-08-13 04:18:49"
$wsDeDe.Range("K2").Value = "2016-08-13 04:19:22"
$wsDeDe.Range("K5").Value = "2016-08-13 04:19:22"
